$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.462.00"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.10%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.422.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.73"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "159.75"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.423.80"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.575"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.34%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.26"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.75%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.435"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.012.24"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.36%  "
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("E15").Value = "  +1.28%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.96"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.52%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.455.72"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.425.83"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("E19").Value = "  -2.54%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.37%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("E22").Value = "  -5.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "72.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.51%  "
$ws.Range("E24").Value = "  -1.45%  "
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  +9.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.25%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("E29").Value = "  +0.32%  "
$ws.Range("E30").Value = "  +3.35%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.41"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.76%  "
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.61%  "
$ws.Range("E34").Value = "  -1.89%  "
$ws.Range("E35").Value = "  +0.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "7.05"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.45%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "162.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.01%  "
$ws.Range("E38").Value = "  -2.21%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.31%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.013.53"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0760"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.55%  "
$ws.Range("E42").Value = "  -5.98%  "
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.63"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.88%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0316"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.70%  "
$ws.Range("E46").Value = "  -2.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "24.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.48%  "
$ws.Range("E48").Value = "  -3.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.867"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.08%  "
$ws.Range("E50").Value = "  +2.50%  "
$ws.Range("B51").Value = "Stellar"
$ws.Range("C51").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.106"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.10%  "
